$d = $word.ActiveDocument
$t = $d.Tables.Item(2)

$t.Cell(3, 6).Range.Text  = "1.801"
$t.Cell(4, 6).Range.Text  = "137.0"
$t.Cell(5, 6).Range.Text  = "1.793"
$t.Cell(6, 6).Range.Text  = "4"
$t.Cell(7, 6).Range.Text  = "2.997"
$t.Cell(8, 6).Range.Text  = "3"
$t.Cell(9, 6).Range.Text  = "5.03"
$t.Cell(10, 6).Range.Text = "1.201"
$t.Cell(11, 6).Range.Text = "0.893"
$t.Cell(12, 6).Range.Text = "0.900"
